$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

Set-TextValue $ws.Range("D2") '37.848.51'
$ws.Range("E2").Value = '  +1.19%  '
Set-TextValue $ws.Range("D3") '2.085.96'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '232.70'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("E7").Value = '  -0.04%  '
Set-TextValue $ws.Range("D8") '57.34'
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("E11").Value = '  +2.80%  '
Set-TextValue $ws.Range("D12") '2.381.88'
$ws.Range("E12").Value = '  +0.50%  '
Set-TextValue $ws.Range("D13") '14.37'
$ws.Range("E13").Value = '  -1.73%  '
Set-TextValue $ws.Range("D14") '21.05'
$ws.Range("E14").Value = '  +1.74%  '
Set-TextValue $ws.Range("D15") '0.761'
$ws.Range("E15").Value = '  -2.15%  '
$ws.Range("E16").Value = '  +1.96%  '
Set-TextValue $ws.Range("D17") '2.077.62'
$ws.Range("E17").Value = '  +0.56%  '
Set-TextValue $ws.Range("D18") '37.753.41'
$ws.Range("E19").Value = '  -2.09%  '
Set-TextValue $ws.Range("D20") '70.88'
$ws.Range("E20").Value = '  +2.07%  '
Set-TextValue $ws.Range("D21") '0.0₃0821'
$ws.Range("E21").Value = '  +1.26%  '
Set-TextValue $ws.Range("D22") '228.05'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -1.87%  '
Set-TextValue $ws.Range("D25") '2.38'
$ws.Range("E25").Value = '  -1.13%  '
Set-TextValue $ws.Range("D26") '170.31'
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("E27").Value = '  +9.40%  '
$ws.Range("E28").Value = '  +1.39%  '
$ws.Range("E29").Value = '  -0.08%  '
Set-TextValue $ws.Range("D30") '19.47'
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("E31").Value = '  +0.56%  '
Set-TextValue $ws.Range("D32") '4.61'
$ws.Range("E32").Value = '  +3.39%  '
Set-TextValue $ws.Range("D33") '0.0625'
$ws.Range("E33").Value = '  +1.21%  '
Set-TextValue $ws.Range("D34") '4.58'
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").Value = '  +3.85%  '
Set-TextValue $ws.Range("D37") '3.38'
$ws.Range("E37").Value = '  +4.43%  '
$ws.Range("E38").Value = '  -0.08%  '
Set-TextValue $ws.Range("D39") '5.41'
$ws.Range("E39").Value = '  -4.83%  '
$ws.Range("E40").Value = '  +5.74%  '
$ws.Range("E41").Value = '  -0.86%  '
Set-TextValue $ws.Range("D42") '97.01'
$ws.Range("E43").Value = '  +0.22%  '
Set-TextValue $ws.Range("D44") '1.451.78'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("E46").Value = '  +3.02%  '
Set-TextValue $ws.Range("D47") '4.06'
$ws.Range("E47").Value = '  -6.94%  '
$ws.Range("E48").Value = '  +3.54%  '
Set-TextValue $ws.Range("D49") '7.37'
$ws.Range("E49").Value = '  +3.12%  '
Set-TextValue $ws.Range("D50") '3.00'
$ws.Range("E50").Value = '  +1.44%  '
Set-TextValue $ws.Range("D51") '2.277.12'
$ws.Range("E51").Value = '  +0.79%  '
